$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Expand the table by one row (this grows the table ref/autoFilter and
# worksheet dimension automatically, matching how Excel grows a Table
# when a new row of data is appended below it)
$table = $ws.ListObjects.Item(1)
$table.ListRows.Add() | Out-Null

# Fill in the missing End Time value for the existing last row (37)
$ws.Range("C37").Value = 0

# New daily power record row (38)
$ws.Range("A38").Value = 43362
$ws.Range("B38").Value = 0.81597222222222221
$ws.Range("C38").Value = 0.99930555555555556
$ws.Range("D38").Formula = "=(C38-B38)* 1440"
$ws.Range("E38").Formula = "=IF(C38>B38, (C38-B38)*1440, (B38-C38)*1440)"
$ws.Range("F38").Formula = "=ABS((C38-B38)*1440)"

# Update the selection/view to the newly-added row
$ws.Range("D38").Select()
$excel.ActiveWindow.ScrollRow = 27
